$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 247-248. This shifts the existing rows 247:260
# down to 249:262 (their contents stay identical, only their row numbers
# change), and leaves two blank rows at 247:248 ready to be populated
# with the new weekly entries.
$ws.Rows("247:248").Insert()

# --- Row 247: Durazno, Doctor Davis, Especial ---
$ws.Range("A247").Value = 7
$ws.Range("B247").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C247").Value = "Ñuble"
$ws.Range("D247").Value = 44615
$ws.Range("E247").Value = 16
$ws.Range("F247").Value = "Fruta"
$ws.Range("G247").Value = 100103
$ws.Range("H247").Value = "Frutos de hueso (carozo)"
$ws.Range("I247").Value = 100103004
$ws.Range("J247").Value = "Durazno"
$ws.Range("K247").Value = "Doctor Davis"
$ws.Range("L247").Value = "Especial"
$ws.Range("M247").Value = 80
$ws.Range("N247").Value = 12000
$ws.Range("O247").Value = 12000
$ws.Range("P247").Value = 12000
$ws.Range("Q247").Value = "$/caja 16 kilos empedrada"
$ws.Range("R247").Value = "Región de O'Higgins"
$ws.Range("S247").Value = 750
$ws.Range("T247").Value = 16

# --- Row 248: Durazno, Doctor Davis, Primera ---
$ws.Range("A248").Value = 7
$ws.Range("B248").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C248").Value = "Ñuble"
$ws.Range("D248").Value = 44615
$ws.Range("E248").Value = 16
$ws.Range("F248").Value = "Fruta"
$ws.Range("G248").Value = 100103
$ws.Range("H248").Value = "Frutos de hueso (carozo)"
$ws.Range("I248").Value = 100103004
$ws.Range("J248").Value = "Durazno"
$ws.Range("K248").Value = "Doctor Davis"
$ws.Range("L248").Value = "Primera"
$ws.Range("M248").Value = 120
$ws.Range("N248").Value = 10000
$ws.Range("O248").Value = 11000
$ws.Range("P248").Value = 10500
$ws.Range("Q248").Value = "$/caja 16 kilos empedrada"
$ws.Range("R248").Value = "Región de O'Higgins"
$ws.Range("S248").Value = 656
$ws.Range("T248").Value = 16
